# Update Scott and Alex's data in "TheBiggestLoser" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TheBiggestLoser")

# Alex's weight (B2) was unknown ("?"); fill in the real number and restore
# the normal (non-highlighted) cell formatting used by its neighbours.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B2").Value = 86.9

# Scott.C's weight (B5) was unknown ("?"); fill in the real number and
# restore the normal cell formatting as well.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B5").Value = 65

# Scott.L's cell (B11) was also highlighted as unknown; now that a concrete
# value is known it gets the regular formatting used by the rest of its team.
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("B11").Value = "X + 3.4"

$excel.CutCopyMode = 0

# Recalculate so the SUM formula in B7 reflects the new figures.
$wb.Application.Calculate() | Out-Null

# Restore the active selection to where the user left off editing.
$ws.Activate() | Out-Null
$ws.Range("D11").Select() | Out-Null
